# Scheduled price/profit refresh for the Leve tables (ALC/ARM/BSM/CRP/CUL/GSM/LTW).
# Updates currentAveragePrice(NQ/HQ) (H:J), LevePrice(NQ/HQ) (K:L) and the
# derived LeveProfit(NQ/HQ) (M:N) columns per row with freshly pulled market data.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 106: Making Your Mark
$ws.Cells.Item(106, 8).Value = 0
$ws.Cells.Item(106, 10).Value = 0
$ws.Cells.Item(106, 12).Value = 0
$ws.Cells.Item(106, 14).ClearContents()
# Row 127: Liquid Competence
$ws.Cells.Item(127, 8).Value = 2167
$ws.Cells.Item(127, 9).Value = 2308.2222
$ws.Cells.Item(127, 11).Value = 6924.6666
$ws.Cells.Item(127, 13).Value = -1964.6666
# Row 129: Practical Command
$ws.Cells.Item(129, 8).Value = 1832
$ws.Cells.Item(129, 9).Value = 1920.3334
$ws.Cells.Item(129, 10).Value = 1699.5
$ws.Cells.Item(129, 11).Value = 5761.0002
$ws.Cells.Item(129, 12).Value = 5098.5
$ws.Cells.Item(129, 13).Value = -761.0002000000004
$ws.Cells.Item(129, 14).Value = -15098.5
# Row 138: All-night Crafting
$ws.Cells.Item(138, 8).Value = 3897.0264
$ws.Cells.Item(138, 10).Value = 4045.375
$ws.Cells.Item(138, 12).Value = 12136.125
$ws.Cells.Item(138, 14).Value = -22416.125

$ws = $wb.Worksheets.Item("ARM")
# Row 32: Ingot We Trust
$ws.Cells.Item(32, 8).Value = 2486022.5
$ws.Cells.Item(32, 9).Value = 2335556.8
$ws.Cells.Item(32, 11).Value = 2335556.8
$ws.Cells.Item(32, 13).Value = -2335269.8
# Row 45: Hollow Hallmarks
$ws.Cells.Item(45, 8).Value = 1791.2632
$ws.Cells.Item(45, 9).Value = 1232.4286
$ws.Cells.Item(45, 10).Value = 3356
$ws.Cells.Item(45, 11).Value = 1232.4286
$ws.Cells.Item(45, 12).Value = 3356
$ws.Cells.Item(45, 13).Value = -855.4286
$ws.Cells.Item(45, 14).Value = -4110
# Row 110: Scheduled Maintenance
$ws.Cells.Item(110, 8).Value = 1232.5294
$ws.Cells.Item(110, 9).Value = 1243.5333
$ws.Cells.Item(110, 11).Value = 1243.5333
$ws.Cells.Item(110, 13).Value = 801.4667
# Row 124: Ace of Gloves
$ws.Cells.Item(124, 8).Value = 54490
$ws.Cells.Item(124, 10).Value = 54490
$ws.Cells.Item(124, 12).Value = 54490
$ws.Cells.Item(124, 14).Value = -64310
# Row 132: Don't Bore Me, Ore Me
$ws.Cells.Item(132, 8).Value = 1776.2667
$ws.Cells.Item(132, 9).Value = 1615
$ws.Cells.Item(132, 11).Value = 4845
$ws.Cells.Item(132, 13).Value = -2315

$ws = $wb.Worksheets.Item("BSM")
# Row 86: Through Thick and Thin
$ws.Cells.Item(86, 8).Value = 7084.1665
$ws.Cells.Item(86, 9).Value = 7801
$ws.Cells.Item(86, 10).Value = 3500
$ws.Cells.Item(86, 11).Value = 7801
$ws.Cells.Item(86, 12).Value = 3500
$ws.Cells.Item(86, 13).Value = -6678
$ws.Cells.Item(86, 14).Value = -5746
# Row 89: Piercing Eyes Deserve Piercing Shafts (L)
$ws.Cells.Item(89, 8).Value = 7084.1665
$ws.Cells.Item(89, 9).Value = 7801
$ws.Cells.Item(89, 10).Value = 3500
$ws.Cells.Item(89, 11).Value = 39005
$ws.Cells.Item(89, 12).Value = 17500
$ws.Cells.Item(89, 13).Value = -33389
$ws.Cells.Item(89, 14).Value = -28732
# Row 107: The Gold Experience
$ws.Cells.Item(107, 8).Value = 5749.25
$ws.Cells.Item(107, 9).Value = 5549
$ws.Cells.Item(107, 11).Value = 5549
$ws.Cells.Item(107, 13).Value = -3629

$ws = $wb.Worksheets.Item("CRP")
# Row 7: Gridania's Got Talent
$ws.Cells.Item(7, 8).Value = 133.75
$ws.Cells.Item(7, 10).Value = 0
$ws.Cells.Item(7, 12).Value = 0
$ws.Cells.Item(7, 14).ClearContents()
# Row 31: Wall Not Found
$ws.Cells.Item(31, 8).Value = 1133.4
$ws.Cells.Item(31, 9).Value = 988
$ws.Cells.Item(31, 10).Value = 1211.6923
$ws.Cells.Item(31, 11).Value = 988
$ws.Cells.Item(31, 12).Value = 1211.6923
$ws.Cells.Item(31, 13).Value = -693
$ws.Cells.Item(31, 14).Value = -1801.6923
# Row 34: Armoires of the Rich and Famous
$ws.Cells.Item(34, 8).Value = 1133.4
$ws.Cells.Item(34, 9).Value = 988
$ws.Cells.Item(34, 10).Value = 1211.6923
$ws.Cells.Item(34, 11).Value = 988
$ws.Cells.Item(34, 12).Value = 1211.6923
$ws.Cells.Item(34, 13).Value = -786
$ws.Cells.Item(34, 14).Value = -1615.6923
# Row 58: You Do the Heavy Lifting
$ws.Cells.Item(58, 8).Value = 3134.2666
$ws.Cells.Item(58, 9).Value = 1973.75
$ws.Cells.Item(58, 11).Value = 1973.75
$ws.Cells.Item(58, 13).Value = -1770.75
# Row 124: Earring Awakening
$ws.Cells.Item(124, 8).Value = 89763.336
$ws.Cells.Item(124, 10).Value = 89763.336
$ws.Cells.Item(124, 12).Value = 89763.336
$ws.Cells.Item(124, 14).Value = -94673.336
# Row 132: Hull Lotta Damage
$ws.Cells.Item(132, 8).Value = 5912.4287
$ws.Cells.Item(132, 9).Value = 5912.4287
$ws.Cells.Item(132, 11).Value = 17737.2861
$ws.Cells.Item(132, 13).Value = -15207.2861
# Row 136: Turali Quality
$ws.Cells.Item(136, 8).Value = 3134.2666
$ws.Cells.Item(136, 9).Value = 1973.75
$ws.Cells.Item(136, 11).Value = 5921.25
$ws.Cells.Item(136, 13).Value = -3371.25

$ws = $wb.Worksheets.Item("CUL")
# Row 5: What a Sap
$ws.Cells.Item(5, 8).Value = 840.2143
$ws.Cells.Item(5, 9).Value = 881.8
$ws.Cells.Item(5, 10).Value = 736.25
$ws.Cells.Item(5, 11).Value = 2645.4
$ws.Cells.Item(5, 12).Value = 2208.75
$ws.Cells.Item(5, 13).Value = -2533.4
$ws.Cells.Item(5, 14).Value = -2432.75
# Row 23: Sweet Smell of Success
$ws.Cells.Item(23, 8).Value = 1596.75
$ws.Cells.Item(23, 9).Value = 450
$ws.Cells.Item(23, 11).Value = 1350
$ws.Cells.Item(23, 13).Value = -1115
# Row 80: Saucy for a Suitor
$ws.Cells.Item(80, 8).Value = 2496.5
$ws.Cells.Item(80, 10).Value = 2991
$ws.Cells.Item(80, 12).Value = 8973
$ws.Cells.Item(80, 14).Value = -10845
# Row 83: Saved by the Sauce (L)
$ws.Cells.Item(83, 8).Value = 2496.5
$ws.Cells.Item(83, 10).Value = 2991
$ws.Cells.Item(83, 12).Value = 26919
$ws.Cells.Item(83, 14).Value = -36279
# Row 131: The Mountain Steeped
$ws.Cells.Item(131, 8).Value = 501775.5
$ws.Cells.Item(131, 9).Value = 1097.5
$ws.Cells.Item(131, 10).Value = 716351.8
$ws.Cells.Item(131, 11).Value = 3292.5
$ws.Cells.Item(131, 12).Value = 2149055.4
$ws.Cells.Item(131, 13).Value = 1747.5
$ws.Cells.Item(131, 14).Value = -2159135.4
# Row 135: Not-so-secret Ingredient
$ws.Cells.Item(135, 8).Value = 840.2143
$ws.Cells.Item(135, 9).Value = 881.8
$ws.Cells.Item(135, 10).Value = 736.25
$ws.Cells.Item(135, 11).Value = 7936.2
$ws.Cells.Item(135, 12).Value = 6626.25
$ws.Cells.Item(135, 13).Value = -5401.2
$ws.Cells.Item(135, 14).Value = -11696.25
# Row 138: Bring Me Your Tacos
$ws.Cells.Item(138, 8).Value = 846
$ws.Cells.Item(138, 9).Value = 846
$ws.Cells.Item(138, 11).Value = 2538
$ws.Cells.Item(138, 13).Value = 2602
# Row 140: Sweet, Sweet Bean Juice
$ws.Cells.Item(140, 8).Value = 9773.9375
$ws.Cells.Item(140, 9).Value = 1735.375
$ws.Cells.Item(140, 11).Value = 5206.125
$ws.Cells.Item(140, 13).Value = -26.125

$ws = $wb.Worksheets.Item("GSM")
# Row 102: Put the Metal to the Peddle
$ws.Cells.Item(102, 8).Value = 2206.8462
$ws.Cells.Item(102, 9).Value = 1687.6666
$ws.Cells.Item(102, 11).Value = 1687.6666
$ws.Cells.Item(102, 13).Value = -65.66660000000002
# Row 107: Whetstones for the Workers
$ws.Cells.Item(107, 8).Value = 1000
$ws.Cells.Item(107, 9).Value = 1000
$ws.Cells.Item(107, 11).Value = 1000
$ws.Cells.Item(107, 13).Value = 920
# Row 113: Copious Crystal Cannons
$ws.Cells.Item(113, 8).Value = 1650
$ws.Cells.Item(113, 9).Value = 0
$ws.Cells.Item(113, 10).Value = 1650
$ws.Cells.Item(113, 11).Value = 0
$ws.Cells.Item(113, 12).Value = 1650
$ws.Cells.Item(113, 13).ClearContents()
$ws.Cells.Item(113, 14).Value = -5990
# Row 132: On Board for Lar
$ws.Cells.Item(132, 8).Value = 3749.6667
$ws.Cells.Item(132, 9).Value = 3499.8
$ws.Cells.Item(132, 11).Value = 10499.4
$ws.Cells.Item(132, 13).Value = -7969.400000000001

$ws = $wb.Worksheets.Item("LTW")
# Row 22: Skin off Their Backs
$ws.Cells.Item(22, 8).Value = 1186.1111
$ws.Cells.Item(22, 10).Value = 1375
$ws.Cells.Item(22, 12).Value = 1375
$ws.Cells.Item(22, 14).Value = -1965
# Row 27: Fire and Hide
$ws.Cells.Item(27, 8).Value = 1186.1111
$ws.Cells.Item(27, 10).Value = 1375
$ws.Cells.Item(27, 12).Value = 1375
$ws.Cells.Item(27, 14).Value = -1589
# Row 32: Men Who Scare Up Goats
$ws.Cells.Item(32, 8).Value = 2963
$ws.Cells.Item(32, 9).Value = 2963
$ws.Cells.Item(32, 11).Value = 2963
$ws.Cells.Item(32, 13).Value = -2646

